$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 16335
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 16335
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H111").Value = 1911.5
$ws.Range("I111").Value = 1093.8
$ws.Range("J111").Value = 6000
$ws.Range("K111").Value = 3281.4
$ws.Range("L111").Value = 18000
$ws.Range("M111").Value = -214.3999999999996
$ws.Range("N111").Value = -24134
$ws.Range("H116").Value = 40062.5
$ws.Range("I116").Value = 42625
$ws.Range("J116").Value = 37500
$ws.Range("K116").Value = 42625
$ws.Range("L116").Value = 37500
$ws.Range("M116").Value = -39183
$ws.Range("N116").Value = -44384
$ws.Range("H132").Value = 11181.25
$ws.Range("I132").Value = 11611.526
$ws.Range("J132").Value = 3006
$ws.Range("K132").Value = 34834.578
$ws.Range("L132").Value = 9018
$ws.Range("M132").Value = -32304.578
$ws.Range("N132").Value = -14078

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 223443.61
$ws.Range("I32").Value = 254229.97
$ws.Range("K32").Value = 254229.97
$ws.Range("M32").Value = -253942.97
$ws.Range("H45").Value = 1743.3334
$ws.Range("I45").Value = 1317.2222
$ws.Range("K45").Value = 1317.2222
$ws.Range("M45").Value = -940.2221999999999
$ws.Range("H61").Value = 2829.2593
$ws.Range("I61").Value = 2578.75
$ws.Range("K61").Value = 2578.75
$ws.Range("M61").Value = -2366.75
$ws.Range("H119").Value = 19833.334
$ws.Range("J119").Value = 19833.334
$ws.Range("L119").Value = 19833.334
$ws.Range("N119").Value = -29509.334
$ws.Range("H132").Value = 5123.904
$ws.Range("I132").Value = 3825.879
$ws.Range("J132").Value = 7378.3687
$ws.Range("K132").Value = 11477.637
$ws.Range("L132").Value = 22135.1061
$ws.Range("M132").Value = -8947.636999999999
$ws.Range("N132").Value = -27195.1061
$ws.Range("H136").Value = 2829.2593
$ws.Range("I136").Value = 2578.75
$ws.Range("K136").Value = 7736.25
$ws.Range("M136").Value = -5186.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 5474.625
$ws.Range("I94").Value = 6049.5
$ws.Range("K94").Value = 6049.5
$ws.Range("M94").Value = -5598.5
$ws.Range("H107").Value = 1528.7838
$ws.Range("I107").Value = 1184.8182
$ws.Range("K107").Value = 1184.8182
$ws.Range("M107").Value = 735.1818000000001
$ws.Range("H134").Value = 4972.7144
$ws.Range("I134").Value = 4921.8
$ws.Range("K134").Value = 14765.4
$ws.Range("M134").Value = -12230.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3316
$ws.Range("H34").Value = 3316
$ws.Range("H59").Value = 25226.25
$ws.Range("J59").Value = 30953.75
$ws.Range("L59").Value = 30953.75
$ws.Range("N59").Value = -33243.75
$ws.Range("H68").Value = 41782.832
$ws.Range("J68").Value = 41782.832
$ws.Range("L68").Value = 41782.832
$ws.Range("N68").Value = -43280.832
$ws.Range("H71").Value = 41782.832
$ws.Range("J71").Value = 41782.832
$ws.Range("L71").Value = 125348.496
$ws.Range("N71").Value = -132836.496
$ws.Range("H134").Value = 2293.9033
$ws.Range("I134").Value = 2055.625
$ws.Range("K134").Value = 6166.875
$ws.Range("M134").Value = -3631.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4309232.5
$ws.Range("I4").Value = 10182664
$ws.Range("K4").Value = 30547992
$ws.Range("M4").Value = -30547880
$ws.Range("H33").Value = 2012
$ws.Range("J33").Value = 2791.4
$ws.Range("L33").Value = 16748.4
$ws.Range("N33").Value = -17314.4
$ws.Range("H134").Value = 7882.9165
$ws.Range("I134").Value = 2000
$ws.Range("K134").Value = 6000
$ws.Range("M134").Value = -930

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 105.681816
$ws.Range("I2").Value = 54.25
$ws.Range("K2").Value = 54.25
$ws.Range("M2").Value = 58.75
$ws.Range("H113").Value = 2126.2173
$ws.Range("I113").Value = 2052.7896
$ws.Range("J113").Value = 2475
$ws.Range("K113").Value = 2052.7896
$ws.Range("L113").Value = 2475
$ws.Range("M113").Value = 117.2103999999999
$ws.Range("N113").Value = -6815
$ws.Range("H126").Value = 2730.5557
$ws.Range("I126").Value = 2735
$ws.Range("J126").Value = 2725
$ws.Range("K126").Value = 8205
$ws.Range("L126").Value = 8175
$ws.Range("M126").Value = -5735
$ws.Range("N126").Value = -13115

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 10234.5
$ws.Range("I93").Value = 500
$ws.Range("J93").Value = 19969
$ws.Range("K93").Value = 500
$ws.Range("L93").Value = 19969
$ws.Range("M93").Value = 748
$ws.Range("N93").Value = -22465
$ws.Range("H132").Value = 2972.0881
$ws.Range("J132").Value = 3678.2666
$ws.Range("L132").Value = 11034.7998
$ws.Range("N132").Value = -16094.7998
$ws.Range("H136").Value = 3035.5
$ws.Range("I136").Value = 3035.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9106.5
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -6556.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1960.6
$ws.Range("I100").Value = 900
$ws.Range("J100").Value = 2667.6667
$ws.Range("K100").Value = 1800
$ws.Range("L100").Value = 5335.3334
$ws.Range("M100").Value = -1259
$ws.Range("N100").Value = -6417.3334
$ws.Range("H107").Value = 1603.5358
$ws.Range("I107").Value = 1137.0625
$ws.Range("J107").Value = 2225.5
$ws.Range("K107").Value = 3411.1875
$ws.Range("L107").Value = 6676.5
$ws.Range("M107").Value = -1491.1875
$ws.Range("N107").Value = -10516.5
$ws.Range("H118").Value = 45000
$ws.Range("J118").Value = 45000
$ws.Range("L118").Value = 45000
$ws.Range("N118").Value = -48314
$ws.Range("H126").Value = 1143.7142
$ws.Range("I126").Value = 1143.7142
$ws.Range("K126").Value = 3431.1426
$ws.Range("M126").Value = -961.1425999999997
$ws.Range("H132").Value = 3539.077
$ws.Range("I132").Value = 2982.3333
$ws.Range("K132").Value = 8946.999899999999
$ws.Range("M132").Value = -6416.999899999999
